# Commit: "Taking latest changes and appending my changes"
#
# - Clear out the stale "Results" column (E2:E41) on the "Test Cases" sheet;
#   the values that used to live there (PASS/SKIP/FAIL) are no longer
#   accurate after merging in the latest changes.
# - Move the active selection from the old Runmode column (D2:D41) to the
#   now-empty Results column (E2:E41), matching where the author was about
#   to fill in fresh results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

# Clear the "Results" column (E) for every data row (2 through 41).
$ws.Range("E2:E41").ClearContents()

# Update the selection to reflect the column the author is now working in.
$ws.Range("E2:E41").Select()
